$d = $word.ActiveDocument

$d.Content.Find.Execute("36×30=1080", $true, $false, $false, $false, $false, $true, 1, $false, "41×43=1763", 2) | Out-Null
$d.Content.Find.Execute("16×43=688", $true, $false, $false, $false, $false, $true, 1, $false, "43×47=2021", 2) | Out-Null
$d.Content.Find.Execute("51×69=3519", $true, $false, $false, $false, $false, $true, 1, $false, "85×80=6800", 2) | Out-Null
$d.Content.Find.Execute("29×87=2523", $true, $false, $false, $false, $false, $true, 1, $false, "66×21=1386", 2) | Out-Null
$d.Content.Find.Execute("44×67=2948", $true, $false, $false, $false, $false, $true, 1, $false, "75×97=7275", 2) | Out-Null
$d.Content.Find.Execute("56×73=4088", $true, $false, $false, $false, $false, $true, 1, $false, "42×82=3444", 2) | Out-Null
$d.Content.Find.Execute("88×24=2112", $true, $false, $false, $false, $false, $true, 1, $false, "31×57=1767", 2) | Out-Null
$d.Content.Find.Execute("12×86=1032", $true, $false, $false, $false, $false, $true, 1, $false, "26×39=1014", 2) | Out-Null
$d.Content.Find.Execute("33×32=1056", $true, $false, $false, $false, $false, $true, 1, $false, "33×36=1188", 2) | Out-Null
$d.Content.Find.Execute("61×53=3233", $true, $false, $false, $false, $false, $true, 1, $false, "86×92=7912", 2) | Out-Null
$d.Content.Find.Execute("71×64=4544", $true, $false, $false, $false, $false, $true, 1, $false, "38×26=988", 2) | Out-Null
$d.Content.Find.Execute("73×13=949", $true, $false, $false, $false, $false, $true, 1, $false, "17×91=1547", 2) | Out-Null
$d.Content.Find.Execute("67×70=4690", $true, $false, $false, $false, $false, $true, 1, $false, "22×97=2134", 2) | Out-Null
$d.Content.Find.Execute("71×86=6106", $true, $false, $false, $false, $false, $true, 1, $false, "21×52=1092", 2) | Out-Null
$d.Content.Find.Execute("58×76=4408", $true, $false, $false, $false, $false, $true, 1, $false, "18×17=306", 2) | Out-Null
$d.Content.Find.Execute("26×72=1872", $true, $false, $false, $false, $false, $true, 1, $false, "99×44=4356", 2) | Out-Null
$d.Content.Find.Execute("88×37=3256", $true, $false, $false, $false, $false, $true, 1, $false, "40×51=2040", 2) | Out-Null
$d.Content.Find.Execute("91×45=4095", $true, $false, $false, $false, $false, $true, 1, $false, "13×68=884", 2) | Out-Null
$d.Content.Find.Execute("16×45=720", $true, $false, $false, $false, $false, $true, 1, $false, "28×39=1092", 2) | Out-Null
$d.Content.Find.Execute("64×62=3968", $true, $false, $false, $false, $false, $true, 1, $false, "62×40=2480", 2) | Out-Null
$d.Content.Find.Execute("21×69=1449", $true, $false, $false, $false, $false, $true, 1, $false, "59×56=3304", 2) | Out-Null
$d.Content.Find.Execute("60×22=1320", $true, $false, $false, $false, $false, $true, 1, $false, "91×53=4823", 2) | Out-Null
$d.Content.Find.Execute("55×19=1045", $true, $false, $false, $false, $false, $true, 1, $false, "30×80=2400", 2) | Out-Null
$d.Content.Find.Execute("23×76=1748", $true, $false, $false, $false, $false, $true, 1, $false, "42×70=2940", 2) | Out-Null
$d.Content.Find.Execute("81×25=2025", $true, $false, $false, $false, $false, $true, 1, $false, "49×12=588", 2) | Out-Null
